$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to be bumped by one
# day (45205 -> 45206) for every data row (rows 2 through 200).
for ($r = 2; $r -le 200; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value = 45206
    }
}
